$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 71: politeness_score (column B) changes from text "3" to the numeric value 3
$ws.Cells.Item(71, 2).Value = 3

# Row 72: new annotation row appended, following the pattern of the old row 71
# but with its own data. politeness_score here stays a text "3" (matching the
# target diff, which keeps it as an inline string), so force text formatting
# before assigning the value.
$ws.Cells.Item(72, 1).Value = "Ruilin"
$ws.Cells.Item(72, 2).NumberFormat = "@"
$ws.Cells.Item(72, 2).Value = "3"
$ws.Cells.Item(72, 2).Style = "Normal"
$ws.Cells.Item(72, 3).Value = "无"
$ws.Cells.Item(72, 4).Value = "QSN"
$ws.Cells.Item(72, 5).Value = "OTH"
$ws.Cells.Item(72, 6).Value = "6649e081-efd7-424b-ac96-c0694d3eea45"
$ws.Cells.Item(72, 7).Value = "HyydRMZC-_annotated.xlsx"
$ws.Cells.Item(72, 8).Value = "In particular, what is exactly Opt attack?"
